# redo game logic and tool upgrading
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the old "Unlock X Upgrades" tiers into a proper "Upgrade to X Tools" progression
# and chain their prerequisites off one another (Wooden -> Copper -> Silver -> Golden).
$ws.Range("A14").Value = "Upgrade to Wooden Tools"

$ws.Range("A15").Value = "Upgrade to Copper Tools"
$ws.Range("F15").Value = "Upgrade to Wooden Tools"

$ws.Range("A16").Value = "Upgrade to Silver Tools"
$ws.Range("F16").Value = "Upgrade to Copper Tools"

$ws.Range("A17").Value = "Upgrade to Golden Tools"
$ws.Range("F17").Value = "Upgrade to Silver Tools"

# Stone cost for "Increase Stockpile" bumped from 10 -> 200
$ws.Range("H5").Value = '{"resourceType":"stone","amount":200}'

# Update the current selection/active cell as recorded in the saved view state.
$ws.Range("E18").Select()
